$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each tuple: row, Coin, Link, Price, Volume(1h)
# Values are prefixed with a leading apostrophe to force text entry
# (matches the source data, which stores these columns as strings,
# e.g. "42.717.16" or "1.00" must not become numbers).
$data = @(
    ,@(2, 'Bitcoin', 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc', '42.717.16', '  -6.49%  ')
    ,@(3, 'Ethereum', 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth', '2.226.44', '  -7.18%  ')
    ,@(4, 'TetherUSD', 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt', '1.00', '  +0.18%  ')
    ,@(5, 'BNB', 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb', '313.98', '  -1.52%  ')
    ,@(6, 'Solana', 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol', '99.41', '  -12.83%  ')
    ,@(7, 'XRP', 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp', '0.568', '  -10.53%  ')
    ,@(8, 'USDC', 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc', '1.00', '  +0.00%  ')
    ,@(9, 'Cardano', 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada', '0.559', '  -10.71%  ')
    ,@(10, 'Avalanche', 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax', '37.14', '  -11.73%  ')
    ,@(11, 'OKB', 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb', '53.92', '  -4.09%  ')
    ,@(12, 'Dogecoin', 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge', '0.0836', '  -10.03%  ')
    ,@(13, 'Polkadot', 'https://coinranking.com/coin/25W7FG7om+polkadot-dot', '7.61', '  -12.94%  ')
    ,@(14, 'TRON', 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx', '0.105', '  -4.49%  ')
    ,@(15, 'Polygon', 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic', '0.872', '  -13.30%  ')
    ,@(16, 'WrappedliquidstakedEther2.0', 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth', '2.561.57', '  -7.24%  ')
    ,@(17, 'Chainlink', 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link', '13.91', '  -12.29%  ')
    ,@(18, 'WrappedEther', 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth', '2.215.52', '  -7.85%  ')
    ,@(19, 'WrappedBTC', 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc', '42.557.55', '  -6.71%  ')
    ,@(20, 'InternetComputer(DFINITY)', 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp', '14.13', '  +5.30%  ')
    ,@(21, 'Uniswap', 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni', '6.61', '  -11.75%  ')
    ,@(22, 'ShibaInu', 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib', '0.0₃0945', '  -12.77%  ')
    ,@(23, 'PancakeSwap', 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake', '3.28', '  -7.51%  ')
    ,@(24, 'Litecoin', 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc', '64.77', '  -13.35%  ')
    ,@(25, 'BitcoinCash', 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch', '234.89', '  -11.16%  ')
    ,@(26, 'ImmutableX', 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx', '2.11', '  -10.44%  ')
    ,@(27, 'Dai', 'https://coinranking.com/coin/MoTuySvg7+dai-dai', '1.00', '  +0.12%  ')
    ,@(28, 'Cosmos', 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom', '10.16', '  -10.44%  ')
    ,@(29, 'Filecoin', 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil', '6.57', '  -13.95%  ')
    ,@(30, 'Toncoin', 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton', '2.17', '  -7.79%  ')
    ,@(31, 'Hedera', 'https://coinranking.com/coin/jad286TjB+hedera-hbar', '0.0882', '  -8.81%  ')
    ,@(32, 'EthereumClassic', 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc', '20.47', '  -10.05%  ')
    ,@(33, 'Monero', 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr', '160.06', '  -7.26%  ')
    ,@(34, 'InjectiveProtocol', 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj', '33.18', '  -15.64%  ')
    ,@(35, 'WEMIXToken', 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix', '2.70', '  -8.54%  ')
    ,@(36, 'LidoDAOToken', 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo', '3.10', '  +1.28%  ')
    ,@(37, 'Stellar', 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm', '0.122', '  -7.52%  ')
    ,@(38, 'ARBITRUM', 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb', '1.87', '  +4.92%  ')
    ,@(39, 'RenderToken', 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr', '4.43', '  -10.11%  ')
    ,@(40, 'Kaspa', 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas', '0.103', '  -12.10%  ')
    ,@(41, 'NEARProtocol', 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near', '3.58', '  -13.18%  ')
    ,@(42, 'VeChain', 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet', '0.0322', '  -11.16%  ')
    ,@(43, 'FirstDigitalUSD', 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd', '1.00', '  +0.34%  ')
    ,@(44, 'Maker', 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr', '1.811.22', '  +8.50%  ')
    ,@(45, 'BitcoinSV', 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv', '89.94', '  -11.05%  ')
    ,@(46, 'Celestia', 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia', '12.04', '  -11.42%  ')
    ,@(47, 'Algorand', 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo', '0.207', '  -14.40%  ')
    ,@(48, 'ordi', 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi', '77.52', '  -11.58%  ')
    ,@(49, 'THORChain', 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune', '5.38', '  -6.05%  ')
    ,@(50, 'MultiversX', 'https://coinranking.com/coin/omwkOTglq+multiversx-egld', '60.88', '  -15.73%  ')
    ,@(51, 'FraxShare', 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs', '8.61', '  -9.23%  ')
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = "'" + $row[1]
    $ws.Cells.Item($r, 3).Value = "'" + $row[2]
    $ws.Cells.Item($r, 4).Value = "'" + $row[3]
    $ws.Cells.Item($r, 5).Value = "'" + $row[4]
}

Write-Output "Done updating rows"
